$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '248.86'
$c.Style = "Normal"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '21.95'
$c.Style = "Normal"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '5.501'
$c.Style = "Normal"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '0.05643'
$c.Style = "Normal"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '3.390'
$c.Style = "Normal"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '6.465'
$c.Style = "Normal"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.8026'
$c.Style = "Normal"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.039'
$c.Style = "Normal"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.1436'
$c.Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07323'
$c.Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 12
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.03158'
$c.Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.02952'
$c.Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.09254'
$c.Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.001673'
$c.Style = "Normal"
$ws.Range("E15").Value = '14BitForexTokenBF'
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '3.276'
$c.Style = "Normal"
$ws.Range("E16").Value = '15MCDexMCB'
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.04753'
$c.Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.0005840'
$c.Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.006468'
$c.Style = "Normal"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '0.005031'
$c.Style = "Normal"
$ws.Range("E20").Value = '19HotbitTokenHTBBestin24h'
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.001052'
$c.Style = "Normal"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.0001505'
$c.Style = "Normal"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.0003211'
$c.Style = "Normal"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '4.071'
$c.Style = "Normal"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.109'
$c.Style = "Normal"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 26
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '0.1275'
$c.Style = "Normal"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 28
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 29
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 30
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 31
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 32
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 33
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 34
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 35
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 36
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 37
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 38
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 39
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.04103'
$c.Style = "Normal"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.1043'
$c.Style = "Normal"
$ws.Range("E41").Value = '40BKEXTokenBKK'
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.002981'
$c.Style = "Normal"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.006915'
$c.Style = "Normal"
$ws.Range("E43").Value = '42KickTokenKICK'
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.008940'
$c.Style = "Normal"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.00005660'
$c.Style = "Normal"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.00000000753'
$c.Style = "Normal"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.7879'
$c.Style = "Normal"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.01654'
$c.Style = "Normal"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.00002107'
$c.Style = "Normal"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.01014'
$c.Style = "Normal"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"

# Row 51
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"
